# Add descriptive city lists under five empires/principalities that
# previously had only a bare name in column H (Си Ся, Полоцкое,
# Переяславское, Черниговское, Новгород-Северское). Each empire keeps its
# name in column H on the row of its first listed city, with subsequent
# cities listed in column B only - matching the existing pattern used for
# every other empire/kingdom further up the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: 16 new rows, starting at the old row 363 ("Империя Си Ся").
$ws.Rows("363:378").Insert()

# The Insert() above pushed the five original single-cell rows (the bare
# "Империя Си Ся" / "Полоцкое Княжество" / "Переяславское Княжество" /
# "Черниговское Княжество" / "Новгород-Северское Княжество" labels that used
# to live at H363:H367) down to H379:H383, right on top of where the new
# city rows go below. Clear them now; their names get re-placed at the
# correct new rows (H363/H367/H371/H374/H380) further down.
$ws.Range("H379:H383").ClearContents()

# -- Империя Си Ся --------------------------------------------------------
$ws.Range("B363").Value = "Иньчуань"
$ws.Range("H363").Value = "Империя Си Ся"
$ws.Range("B364").Value = "Цзинань"
$ws.Range("B365").Value = "Сиань"
$ws.Range("B366").Value = "Кайфын"

# -- Полоцкое Княжество ----------------------------------------------------
$ws.Range("B367").Value = "Полоцк"
$ws.Range("H367").Value = "Полоцкое Княжество"
$ws.Range("B368").Value = "Герцике"
$ws.Range("B369").Value = "Борисов"
$ws.Range("B370").Value = "Друцк"

# -- Переяславское Княжество -------------------------------------------------
$ws.Range("B371").Value = "Прилуки"
$ws.Range("H371").Value = "Переяславское Княжество"
$ws.Range("B372").Value = "Лубно"
$ws.Range("B373").Value = "Лукомль"

# -- Черниговское Княжество --------------------------------------------------
$ws.Range("B374").Value = "Меценск"
$ws.Range("H374").Value = "Черниговское Княжество"
$ws.Range("B375").Value = "Стародуб"
$ws.Range("B376").Value = "Рыльск"
$ws.Range("B377").Value = "Рогачев"
$ws.Range("B378").Value = "Гомель"
$ws.Range("B379").Value = "Кромы"

# -- Новгород-Северское Княжество --------------------------------------------
$ws.Range("B380").Value = "Путивль"
$ws.Range("H380").Value = "Новгород-Северское Княжество"
$ws.Range("B381").Value = "Севск"
$ws.Range("B382").Value = "Глухов"
$ws.Range("B383").Value = "Ольгов"

# Restore the view/selection the author ended up with.
$ws.Range("B384").Select()
